$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Range("C38").Value = 3.061100172428026
$ws.Range("D38").Value = 15.37973025068279
$ws.Range("E38").Value = 4.337748530745557
$ws.Range("F38").Value = 12.45999248566049
$ws.Range("G38").Value = 32.17747126708883
$ws.Range("H38").Value = 27.13881411214454
$ws.Range("I38").Value = 7.457878298668167
$ws.Range("J38").Value = 0.5042105153070779
$ws.Range("N38").Value = 0.6279803126980887
$ws.Range("P38").Value = 2.419221143723876
$ws.Range("Q38").Value = 3.061100172428026
$ws.Range("R38").Value = 62.75637990999006
$ws.Range("S38").Value = 0.4532648815865458
$ws.Range("T38").Value = 0.5643759926976569
$ws.Range("V38").Value = -0.1911040132728621
$ws.Range("W38").Value = 32.80884103438532
$ws.Range("X38").Value = 19.71747878142834
$ws.Range("AA38").Value = 12.79379381143391
$ws.Range("C39").Value = 3.098976908803056
$ws.Range("D39").Value = 15.57003241507165
$ws.Range("E39").Value = 4.391421964578305
$ws.Range("F39").Value = 12.29898522126248
$ws.Range("G39").Value = 32.26043960091244
$ws.Range("H39").Value = 26.66527476125871
$ws.Range("I39").Value = 8.790906424986197
$ws.Range("J39").Value = 0.5071117505196556
$ws.Range("L39").Value = 1.180312308138719
$ws.Range("N39").Value = 0.783624314468355
$ws.Range("P39").Value = 3.19574158533247
$ws.Range("Q39").Value = 3.098976908803056
$ws.Range("R39").Value = 63.53290035159866
$ws.Range("S39").Value = 0.5656059193635897
$ws.Range("T39").Value = 0.8156059193635897
$ws.Range("V39").Value = -0.153146802360735
$ws.Range("W39").Value = 33.04745336997919
$ws.Range("X39").Value = 19.96145437964996
$ws.Range("Z39").Value = 2.737759423429122
$ws.Range("AA39").Value = 12.75357792856742
$ws.Range("C40").Value = 3.098781806703383
$ws.Range("D40").Value = 15.56905217349337
$ws.Range("E40").Value = 4.391145494094316
$ws.Range("F40").Value = 12.30371514166842
$ws.Range("G40").Value = 32.2639128092561
$ws.Range("H40").Value = 26.66992013483354
$ws.Range("I40").Value = 8.785734422834107
$ws.Range("J40").Value = 0.5072449932699233
$ws.Range("L40").Value = 1.180312308138719
$ws.Range("P40").Value = 3.191741748411549
$ws.Range("Q40").Value = 3.098781806703383
$ws.Range("R40").Value = 63.52890051467773
$ws.Range("T40").Value = 0.8124853349808941
$ws.Range("W40").Value = 33.04660313382924
$ws.Range("X40").Value = 19.96019766758768
$ws.Range("Z40").Value = 2.737759423429122
$ws.Range("AA40").Value = 12.75430801205244
$ws.Range("C41").Value = 3.068123848016254
$ws.Range("D41").Value = 15.41501894750109
$ws.Range("E41").Value = 4.347701468169067
$ws.Range("F41").Value = 12.28965475459361
$ws.Range("G41").Value = 32.05237517026375
$ws.Range("H41").Value = 26.97363187478792
$ws.Range("I41").Value = 7.641958568352877
$ws.Range("J41").Value = 0.4993209750534031
$ws.Range("N41").Value = 0.783624314468355
$ws.Range("P41").Value = 2.563215272877046
$ws.Range("Q41").Value = 3.068123848016254
$ws.Range("R41").Value = 62.90037403914323
$ws.Range("S41").Value = 0.5656059193635897
$ws.Range("T41").Value = 0.6767170304747009
$ws.Range("V41").Value = -0.153146802360735
$ws.Range("W41").Value = 32.83938893933052
$ws.Range("X41").Value = 19.76272041567015
$ws.Range("AA41").Value = 12.76745020952018
